$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AE2").Value = 980
$ws.Range("AH2").Value = 46
$ws.Range("H2").Value = 2.26
$ws.Range("J2").Value = 3.35
$ws.Range("L2").Value = 1.4
$ws.Range("Q2").Value = 1.87
$ws.Range("R2").Value = 1.37
$ws.Range("S2").Value = 3.2
$ws.Range("T2").Value = 1.69
$ws.Range("V2").Value = 1.66
$ws.Range("Z2").Value = 32
$ws.Range("AB3").Value = 6.6
$ws.Range("AH3").Value = 29
$ws.Range("AN3").Value = 19
$ws.Range("F3").Value = 1.77
$ws.Range("G3").Value = 1.81
$ws.Range("H3").Value = 5.9
$ws.Range("N3").Value = 2.92
$ws.Range("P3").Value = 1.62
$ws.Range("Q3").Value = 2.36
$ws.Range("S3").Value = 4.8
$ws.Range("T3").Value = 2.18
$ws.Range("Y3").Value = 16
$ws.Range("Z3").Value = 46
$ws.Range("F4").Value = 3.6
$ws.Range("H4").Value = 1.85
$ws.Range("I4").Value = 2.06
$ws.Range("J4").Value = 3.6
$ws.Range("K4").Value = 4.7
$ws.Range("M4").Value = 1.06
$ws.Range("P4").Value = 1.79
$ws.Range("Q4").Value = 1.72
$ws.Range("V4").Value = 1.94
$ws.Range("W4").Value = 1.25
$ws.Range("Z4").Value = 40
$ws.Range("AA5").Value = 130
$ws.Range("AB5").Value = 12
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 13
$ws.Range("AE5").Value = 75
$ws.Range("AF5").Value = 21
$ws.Range("AG5").Value = 14
$ws.Range("AH5").Value = 18.5
$ws.Range("AK5").Value = 100
$ws.Range("AN5").Value = 500
$ws.Range("AO5").Value = 29
$ws.Range("F5").Value = 2.86
$ws.Range("G5").Value = 3.15
$ws.Range("Q5").Value = 2.04
$ws.Range("R5").Value = 1.31
$ws.Range("T5").Value = 1.78
$ws.Range("W5").Value = 1.46
$ws.Range("X5").Value = 13.5
$ws.Range("Y5").Value = 11
$ws.Range("Z5").Value = 18
$ws.Range("L6").Value = 1.62
$ws.Range("O6").Value = 1.63
$ws.Range("S6").Value = 5.8
$ws.Range("W6").Value = 1.86
$ws.Range("O7").Value = 1.52
$ws.Range("AF8").Value = 15
$ws.Range("G8").Value = 2.38
$ws.Range("H8").Value = 3.9
$ws.Range("J8").Value = 2.92
$ws.Range("O8").Value = 1.62
$ws.Range("Q8").Value = 2.94
$ws.Range("T8").Value = 2.28
$ws.Range("W8").Value = 1.72
$ws.Range("AD9").Value = 34
$ws.Range("AE9").Value = 510
$ws.Range("AF9").Value = 9.800000000000001
$ws.Range("AH9").Value = 25
$ws.Range("AN9").Value = 5.1
$ws.Range("G9").Value = 1.41
$ws.Range("N9").Value = 5.5
$ws.Range("P9").Value = 2.56
$ws.Range("R9").Value = 1.62
$ws.Range("S9").Value = 2.3
$ws.Range("T9").Value = 1.82
$ws.Range("W9").Value = 3.4
$ws.Range("Z9").Value = 95
$ws.Range("AH10").Value = 30
$ws.Range("AO10").Value = 6.2
$ws.Range("F10").Value = 4.7
$ws.Range("G10").Value = 5.2
$ws.Range("H10").Value = 1.67
$ws.Range("I10").Value = 1.74
$ws.Range("J10").Value = 4.5
$ws.Range("K10").Value = 5
$ws.Range("O10").Value = 1.15
$ws.Range("Q10").Value = 1.47
$ws.Range("S10").Value = 2.2
$ws.Range("V10").Value = 2.34
$ws.Range("W10").Value = 1.24
$ws.Range("Z10").Value = 26
$ws.Range("AI11").Value = 44
$ws.Range("I11").Value = 3.05
$ws.Range("N11").Value = 3.75
$ws.Range("V11").Value = 1.49
$ws.Range("W11").Value = 1.62
$ws.Range("AA12").Value = 42
$ws.Range("AC12").Value = 8.199999999999999
$ws.Range("AH12").Value = 25
$ws.Range("AI12").Value = 980
$ws.Range("AK12").Value = 480
$ws.Range("F12").Value = 4.9
$ws.Range("H12").Value = 1.9
$ws.Range("K12").Value = 3.6
$ws.Range("P12").Value = 1.64
$ws.Range("Q12").Value = 2.32
$ws.Range("X12").Value = 11.5
$ws.Range("Z12").Value = 10.5
$ws.Range("AA13").Value = 32
$ws.Range("AE13").Value = 29
$ws.Range("AH13").Value = 21
$ws.Range("AO13").Value = 25
$ws.Range("H13").Value = 2.2
$ws.Range("I13").Value = 2.24
$ws.Range("N13").Value = 3
$ws.Range("V13").Value = 1.8
$ws.Range("Y13").Value = 8.6
$ws.Range("Z13").Value = 13.5
$ws.Range("G14").Value = 2
$ws.Range("W14").Value = 2
$ws.Range("AN15").Value = 9.199999999999999
$ws.Range("G15").Value = 1.94
$ws.Range("AL16").Value = 46
$ws.Range("AM16").Value = 90
$ws.Range("F16").Value = 2.98
$ws.Range("L17").Value = 1.33
$ws.Range("M17").Value = 1.05
$ws.Range("R17").Value = 1.41
$ws.Range("Z17").Value = 420
$ws.Range("AA19").Value = 900
$ws.Range("L19").Value = 1.34
$ws.Range("Q19").Value = 1.81
$ws.Range("AI20").Value = 44
$ws.Range("AJ20").Value = 60
$ws.Range("J20").Value = 3.3
$ws.Range("N20").Value = 3.55
$ws.Range("S20").Value = 3.95
$ws.Range("T20").Value = 1.87
$ws.Range("AA21").Value = 500
$ws.Range("AG21").Value = 24
$ws.Range("AH21").Value = 28
$ws.Range("F21").Value = 3.1
$ws.Range("G21").Value = 3.7
$ws.Range("H21").Value = 2.24
$ws.Range("J21").Value = 3.1
$ws.Range("K21").Value = 3.65
$ws.Range("O21").Value = 1.34
$ws.Range("P21").Value = 1.73
$ws.Range("Q21").Value = 2.02
$ws.Range("R21").Value = 1.28
$ws.Range("S21").Value = 3.55
$ws.Range("V21").Value = 1.63
$ws.Range("W21").Value = 1.37
$ws.Range("X21").Value = 14
$ws.Range("S22").Value = 3.55
$ws.Range("AB23").Value = 16
$ws.Range("AD23").Value = 10
$ws.Range("AO23").Value = 15.5
$ws.Range("I23").Value = 1.88
$ws.Range("V23").Value = 2.12
$ws.Range("F24").Value = 1.93
$ws.Range("N24").Value = 2.42
$ws.Range("P24").Value = 1.47
$ws.Range("T24").Value = 2.5
$ws.Range("U24").Value = 1.64
$ws.Range("AD25").Value = 16
$ws.Range("AJ25").Value = 370
$ws.Range("AM25").Value = 330
$ws.Range("AN25").Value = 90
$ws.Range("H25").Value = 2.94
$ws.Range("K25").Value = 2.96
$ws.Range("Q25").Value = 3.15
$ws.Range("AK26").Value = 32
$ws.Range("AM26").Value = 340
$ws.Range("I26").Value = 5.7
$ws.Range("K26").Value = 3.25
$ws.Range("V26").Value = 1.22
$ws.Range("Z26").Value = 75
$ws.Range("F27").Value = 3.9
$ws.Range("I27").Value = 2.06
$ws.Range("K27").Value = 5.7
$ws.Range("N27").Value = 2.4
$ws.Range("P27").Value = 1.54
$ws.Range("R27").Value = 1.15
$ws.Range("V27").Value = 1.95
$ws.Range("W27").Value = 1.18
